$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the interrogation dialogue lines in column B (rows 2-17) ---
# Shared strings get re-compacted automatically on save, so we just need to
# push the new text into each cell; unused old strings are dropped and the
# new ones appended by the engine.

$ws.Range("B2").Value  = "I have nothing to do with the Lord’s death!"
$ws.Range("B3").Value  = "That’s not for you to decide. As you saw just now, Ming has entrusted me to investigate the matter."
$ws.Range("B4").Value  = "Hmph. Fine——ask away."
$ws.Range("B5").Value  = "Why are you staying at the manor?"
$ws.Range("B6").Value  = "About half a month ago, I met Ming in town at the foot of the mountain. I complimented his calligraphy——it wasn’t bad."
$ws.Range("B7").Value  = "Then he insisted on inviting me to stay at the manor, so I figured, why not?"
$ws.Range("B8").Value  = "When was the last time you saw the Lord?"
$ws.Range("B9").Value  = "I don’t quite remember......maybe during lunch?"
$ws.Range("B10").Value = "To be honest, I hardly ever leave my room except for meals."
$ws.Range("B11").Value = "Where were you before and after the evening banquet? What were you doing?"
$ws.Range("B12").Value = "Before the banquet started, I wandered around the manor a bit."
$ws.Range("B13").Value = "I recall that after it was confirmed the banquet would be delayed, you left the hall. Where did you go?"
$ws.Range("B14").Value = "I didn’t go anywhere——I went back to my room and slept."
$ws.Range("B15").Value = "I only came here because that butler banged on my door and insisted I join this gathering."
$ws.Range("B16").Value = "So after returning to your room, you never left again?"
$ws.Range("B17").Value = "Exactly! Aren’t you done yet? I already told you——I was asleep the whole time! I didn’t see or hear a thing!"

# Row 14 now holds a longer line that wraps onto two lines in the real sheet.
$ws.Rows(14).RowHeight = 34

# --- Drop the now-empty trailing helper row ---
$ws.Rows(19).Delete()

# --- Restore the selection the author left the sheet in ---
$ws.Range("B19").Select()
